$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 195 ("「あなたであれ。ありのままが美しいから」" post), shifting subsequent rows up
$ws.Rows.Item(195).Delete()
